$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new blank row at position 13 (pushes old rows 13-24 down to
#    14-25, carrying their original row heights/styles with them).
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# The insert leaves a stray empty, styled cell behind in A13; the target
# layout has no cell at all in column A for row 13, so clear it fully.
$ws.Range("A13").Clear()

# ---------------------------------------------------------------------------
# 2. Row 13 (new): "Docentes responsaveis:" value, in B13/C13.
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = "5817181 - Valdeir Arantes"
$ws.Range("C13").Value = "5817181 - Valdeir Arantes"
# Newly-created cells in this row mis-inherit column A's bold style; fix B13
# by pasting the correct (wrap-text, black) format from the row below.
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3. Row 14 ("Programa resumido:") - replace "Semestral" with the real
#    Portuguese short-syllabus text.
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = "1. O papel da estatística na Engenharia2. Fundamentos de estatística aplicada3. Análise de Variância4. Testes de comparações múltiplas5. Planejamento de Experimentos"
$ws.Range("C14").Value = "1. O papel da estatística na Engenharia2. Fundamentos de estatística aplicada3. Análise de Variância4. Testes de comparações múltiplas5. Planejamento de Experimentos"

# Row 15 ("Short syllabus:") keeps its original content untouched.

# ---------------------------------------------------------------------------
# 4. Row 16 ("Programa:") - fix the old bug where B/C mistakenly reused the
#    "01/01/2019" string; set the real Portuguese syllabus text.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "1. O papel da estatística na Engenharia: métodos de coleta de dados2. Fundamentos de estatística aplicada3. Análise de Variância: análise de variância de um modelo4. Testes de comparações múltiplas (Tukey, Hsu)5. Planejamento de Experimentos: vantagens dos experimentos fatoriais em relação aos experimentos do tipo um fator por vez; varielaboração do planejamento fatorial Completo do tipo 2^k e fracionado, e superfície de resposta"
$ws.Range("C16").Value = "1. O papel da estatística na Engenharia: métodos de coleta de dados2. Fundamentos de estatística aplicada3. Análise de Variância: análise de variância de um modelo4. Testes de comparações múltiplas (Tukey, Hsu)5. Planejamento de Experimentos: vantagens dos experimentos fatoriais em relação aos experimentos do tipo um fator por vez; varielaboração do planejamento fatorial Completo do tipo 2^k e fracionado, e superfície de resposta"

# Row 17 ("Syllabus:") keeps its original content untouched.
# Row 18 ("Avaliação:") keeps its original content untouched.

# ---------------------------------------------------------------------------
# 5. Row 19 ("Método:") - fix the old bug where B/C mistakenly reused the
#    "5817181 - Valdeir Arantes" string; set the real method description.
# ---------------------------------------------------------------------------
$ws.Range("B19").Value = "A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# ---------------------------------------------------------------------------
# 6. Row 20 ("Critério:") - the criteria text moves to row 19; this row now
#    takes what used to be the "Norma de recuperação" value.
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = "MF≥ 5,0 para aprovação 5,0"
$ws.Range("C20").Value = "MF≥ 5,0 para aprovação 5,0"

# ---------------------------------------------------------------------------
# 7. Row 21 ("Norma de recuperação:") - takes what used to be the
#    "Bibliografia" value.
# ---------------------------------------------------------------------------
$ws.Range("B21").Value = "(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada."

# ---------------------------------------------------------------------------
# 8. Row 22 ("Bibliografia:") - set the real bibliography text.
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = "1. BOX, G.E.P.; HUNTER, W.G.; HUNTER, J.S. Statistics for Experimenters: an introduction to designs, data analysis and model building. New York: John Wiley & Sons Inc., 1978.2. RODRIGUES, M. I. e IEMMA, A. F. Planejamento de experimentos e otimização de processos. Campinas: Cárita editora, 2009.3. Planejamento e otimização de Experimentos. Roy E. Bruns, Edit. UNICAMP, 1996"
$ws.Range("C22").Value = "1. BOX, G.E.P.; HUNTER, W.G.; HUNTER, J.S. Statistics for Experimenters: an introduction to designs, data analysis and model building. New York: John Wiley & Sons Inc., 1978.2. RODRIGUES, M. I. e IEMMA, A. F. Planejamento de experimentos e otimização de processos. Campinas: Cárita editora, 2009.3. Planejamento e otimização de Experimentos. Roy E. Bruns, Edit. UNICAMP, 1996"

# Rows 23 ("Requisitos:"), 24 and 25 (the two weak-requirement lines) keep
# their original content untouched.

# ---------------------------------------------------------------------------
# 9. Column layout fix: column A's width definition used to overlap column
#    B's (min=1 max=2 followed by a min=2 max=2 override). Re-apply column
#    A's own width so the sheet no longer describes column B via two
#    conflicting <col> ranges.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 29.83716
